$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors used by the "Completed?" status highlighting (stored as old-style
# BGR Long values, same as VBA's RGB()/Interior.Color expects).
$greenFg = 6749952   # RGB(0,255,102)  -> fgColor FF00FF66
$cyanBg  = 16776960  # RGB(0,255,255)  -> bgColor FF00FFFF
$redFg   = 255       # RGB(255,0,0)    -> fgColor FFFF0000
$brownBg = 13209      # RGB(153,51,0)   -> bgColor FF993300

function Set-Done($range) {
    $range.Interior.Color = $greenFg
    $range.Interior.PatternColor = $cyanBg
}

function Set-NotDone($range) {
    $range.Interior.Color = $redFg
    $range.Interior.PatternColor = $brownBg
}

# --- "Menu graphics" (row 3) is now completed ---
$ws.Range("D3").Value = "yes"
Set-Done $ws.Range("B3:D3")

# --- Insert a new task row before "Alter code structure for testing" ---
$ws.Rows("13:13").Insert()

$ws.Range("B13").Value = "Display previous actions"
$ws.Range("C13").Value = "1 day"
$ws.Range("D13").Value = "yes"
Set-Done $ws.Range("B13:D13")

# --- "Test code for game mechanics" (now row 15) is now completed ---
$ws.Range("D15").Value = "yes"
Set-Done $ws.Range("B15:D15")

# --- "Documents / User's manual" (now row 17) is now completed ---
$ws.Range("D17").Value = "yes"
Set-Done $ws.Range("B17:D17")

# --- "Maintenance plan" (now row 19) gains an estimate and is completed ---
$ws.Range("C19").Value = "2-3 days"
$ws.Range("D19").Value = "yes"
Set-Done $ws.Range("B19:D19")

$ws.Range("E14").Select()
